# Auto-generated Excel COM-interop script applying the Zodiark_Profits.xlsx commit diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across the 8 class sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 157.8125
$ws.Range("I2").Value = 157.8125
$ws.Range("K2").Value = 157.8125
$ws.Range("M2").Value = -44.8125
$ws.Range("H29").Value = 6670
$ws.Range("J29").Value = 12500
$ws.Range("L29").Value = 37500
$ws.Range("N29").Value = -38062
$ws.Range("H38").Value = 7824.4614
$ws.Range("I38").Value = 803.8
$ws.Range("J38").Value = 12212.375
$ws.Range("K38").Value = 2411.4
$ws.Range("L38").Value = 36637.125
$ws.Range("M38").Value = -2039.4
$ws.Range("N38").Value = -37381.125
$ws.Range("H51").Value = 5919.5713
$ws.Range("J51").Value = 5919.5713
$ws.Range("L51").Value = 5919.5713
$ws.Range("N51").Value = -6887.5713
$ws.Range("H58").Value = 7271.857
$ws.Range("J58").Value = 16666.334
$ws.Range("L58").Value = 49999.00199999999
$ws.Range("N58").Value = -50299.00199999999
$ws.Range("H88").Value = 1003468.6
$ws.Range("I88").Value = 3711
$ws.Range("K88").Value = 3711
$ws.Range("M88").Value = -3305
$ws.Range("H91").Value = 1003468.6
$ws.Range("I91").Value = 3711
$ws.Range("K91").Value = 3711
$ws.Range("M91").Value = -2307
$ws.Range("H120").Value = 69000
$ws.Range("J120").Value = 69000
$ws.Range("L120").Value = 69000
$ws.Range("N120").Value = -78676
$ws.Range("H132").Value = 4539.8184
$ws.Range("I132").Value = 2457.9644
$ws.Range("K132").Value = 7373.8932
$ws.Range("M132").Value = -4843.8932
$ws.Range("H137").Value = 2114.5264
$ws.Range("I137").Value = 2579.3635
$ws.Range("J137").Value = 1475.375
$ws.Range("K137").Value = 7738.0905
$ws.Range("L137").Value = 4426.125
$ws.Range("M137").Value = -5188.0905
$ws.Range("N137").Value = -9526.125
$ws.Range("H138").Value = 2236.3408
$ws.Range("I138").Value = 1353.2858
$ws.Range("J138").Value = 2403.4055
$ws.Range("K138").Value = 4059.8574
$ws.Range("L138").Value = 7210.2165
$ws.Range("M138").Value = 1080.1426
$ws.Range("N138").Value = -17490.2165

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2151.8333
$ws.Range("I45").Value = 2227.75
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 2227.75
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1850.75
$ws.Range("N45").Value = -2754
$ws.Range("H61").Value = 5595.1665
$ws.Range("I61").Value = 5142.8
$ws.Range("K61").Value = 5142.8
$ws.Range("M61").Value = -4930.8
$ws.Range("H69").Value = 306414.5
$ws.Range("J69").Value = 306414.5
$ws.Range("L69").Value = 306414.5
$ws.Range("N69").Value = -307912.5
$ws.Range("H72").Value = 306414.5
$ws.Range("J72").Value = 306414.5
$ws.Range("L72").Value = 919243.5
$ws.Range("N72").Value = -926731.5
$ws.Range("H103").Value = 95000
$ws.Range("J103").Value = 95000
$ws.Range("L103").Value = 95000
$ws.Range("N103").Value = -97344
$ws.Range("H110").Value = 1670.3334
$ws.Range("I110").Value = 1601.8
$ws.Range("K110").Value = 1601.8
$ws.Range("M110").Value = 443.2
$ws.Range("H122").Value = 3687.25
$ws.Range("I122").Value = 3333.0667
$ws.Range("K122").Value = 9999.2001
$ws.Range("M122").Value = -7549.2001
$ws.Range("H132").Value = 6180.7393
$ws.Range("I132").Value = 7754.222
$ws.Range("K132").Value = 23262.666
$ws.Range("M132").Value = -20732.666
$ws.Range("H136").Value = 5595.1665
$ws.Range("I136").Value = 5142.8
$ws.Range("K136").Value = 15428.4
$ws.Range("M136").Value = -12878.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2202.4443
$ws.Range("I105").Value = 1546.1428
$ws.Range("K105").Value = 1546.1428
$ws.Range("M105").Value = 200.8571999999999
$ws.Range("H134").Value = 10324.042
$ws.Range("I134").Value = 9806.303
$ws.Range("J134").Value = 11463.066
$ws.Range("K134").Value = 29418.909
$ws.Range("L134").Value = 34389.198
$ws.Range("M134").Value = -26883.909
$ws.Range("N134").Value = -39459.198

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3881.9556
$ws.Range("I31").Value = 1212.5714
$ws.Range("K31").Value = 1212.5714
$ws.Range("M31").Value = -917.5714
$ws.Range("H34").Value = 3881.9556
$ws.Range("I34").Value = 1212.5714
$ws.Range("K34").Value = 1212.5714
$ws.Range("M34").Value = -1010.5714
$ws.Range("H86").Value = 111114216
$ws.Range("J86").Value = 3499.5
$ws.Range("L86").Value = 3499.5
$ws.Range("N86").Value = -5745.5
$ws.Range("H89").Value = 111114216
$ws.Range("J89").Value = 3499.5
$ws.Range("L89").Value = 17497.5
$ws.Range("N89").Value = -28729.5
$ws.Range("H122").Value = 2268.7144
$ws.Range("I122").Value = 2230.1667
$ws.Range("K122").Value = 6690.500100000001
$ws.Range("M122").Value = -4240.500100000001
$ws.Range("H134").Value = 1887.84
$ws.Range("I134").Value = 1758.1666
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 5274.4998
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -2739.4998
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 645215.5
$ws.Range("I2").Value = 952407.7
$ws.Range("J2").Value = 111.9
$ws.Range("K2").Value = 5714446.199999999
$ws.Range("L2").Value = 671.4000000000001
$ws.Range("M2").Value = -5714333.199999999
$ws.Range("N2").Value = -897.4000000000001
$ws.Range("H107").Value = 1122.7142
$ws.Range("J107").Value = 1575.4286
$ws.Range("L107").Value = 4726.2858
$ws.Range("N107").Value = -8566.2858
$ws.Range("H129").Value = 4765956
$ws.Range("I129").Value = 5184.5
$ws.Range("J129").Value = 9093930
$ws.Range("K129").Value = 15553.5
$ws.Range("L129").Value = 27281790
$ws.Range("M129").Value = -10553.5
$ws.Range("N129").Value = -27291790
$ws.Range("H131").Value = 2942712
$ws.Range("J131").Value = 1787327.6
$ws.Range("L131").Value = 5361982.800000001
$ws.Range("N131").Value = -5372062.800000001
$ws.Range("H132").Value = 1845.6364
$ws.Range("I132").Value = 1613.5
$ws.Range("J132").Value = 2124.2
$ws.Range("K132").Value = 14521.5
$ws.Range("L132").Value = 19117.8
$ws.Range("M132").Value = -11991.5
$ws.Range("N132").Value = -24177.8
$ws.Range("H139").Value = 4010.077
$ws.Range("I139").Value = 1577.6
$ws.Range("J139").Value = 7327.091
$ws.Range("K139").Value = 4732.799999999999
$ws.Range("L139").Value = 21981.273
$ws.Range("M139").Value = 407.2000000000007
$ws.Range("N139").Value = -32261.273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 48758470
$ws.Range("I126").Value = 11296
$ws.Range("J126").Value = 195000000
$ws.Range("K126").Value = 33888
$ws.Range("L126").Value = 585000000
$ws.Range("M126").Value = -31418
$ws.Range("N126").Value = -585004940
$ws.Range("H132").Value = 7906.9453
$ws.Range("I132").Value = 7553.091
$ws.Range("J132").Value = 8988.166999999999
$ws.Range("K132").Value = 22659.273
$ws.Range("L132").Value = 26964.501
$ws.Range("M132").Value = -20129.273
$ws.Range("N132").Value = -32024.501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1574.7778
$ws.Range("I16").Value = 1351.1538
$ws.Range("K16").Value = 1351.1538
$ws.Range("M16").Value = -1181.1538
$ws.Range("H100").Value = 2736.3333
$ws.Range("J100").Value = 4047
$ws.Range("L100").Value = 4047
$ws.Range("N100").Value = -5129
$ws.Range("H132").Value = 5094.846
$ws.Range("I132").Value = 4964.8647
$ws.Range("J132").Value = 7499.5
$ws.Range("K132").Value = 14894.5941
$ws.Range("L132").Value = 22498.5
$ws.Range("M132").Value = -12364.5941
$ws.Range("N132").Value = -27558.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4223.5
$ws.Range("I122").Value = 4147.175
$ws.Range("K122").Value = 12441.525
$ws.Range("M122").Value = -9991.525000000001
$ws.Range("H132").Value = 3984.3914
$ws.Range("I132").Value = 3438.8386
$ws.Range("J132").Value = 5111.8667
$ws.Range("K132").Value = 10316.5158
$ws.Range("L132").Value = 15335.6001
$ws.Range("M132").Value = -7786.515800000001
$ws.Range("N132").Value = -20395.6001
$ws.Range("H140").Value = 84466.664
$ws.Range("J140").Value = 88000
$ws.Range("L140").Value = 88000
$ws.Range("N140").Value = -98360
